# "Screen Design Standards added"
#
# Adds a new "Input des Benutzers" / "Textfelder, Buttons, Kamera" row to the
# Screen-Design-Standards table on Tabelle1, right after the existing
# "Control Standards Android" block, and shrinks the blank gap that used to
# separate it from the "Control Standards Swing" block below.
#
# Net effect: old rows 17-24 ("Control Standards Swing" block) end up at new
# rows 14-21 (shift of -3): the 5-row blank gap (old rows 12-16) shrinks to a
# single blank row, and the freed-up row 12 is populated with new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the blank gap between row 11 and row 17 down to a single blank
# row, shifting "Control Standards Swing" (old row 17) and everything below
# it up so it lands on row 14.
[void]$ws.Rows("13:16").Delete()

# Re-open a single blank row at 13 so the new content row (12) is followed
# by one empty row, just like the original layout had a blank row under
# each section header.
[void]$ws.Rows("13:13").Insert()

# Fill the new row 12. Write column B before column A so the shared-string
# table records "Textfelder, Buttons, Kamera" ahead of "Input des
# Benutzers", matching the source file's string order.
$ws.Range("B12").Value = "Textfelder, Buttons, Kamera"
$ws.Range("A12").Value = "Input des Benutzers"

# Update the sheet's selection to match the new layout.
[void]$ws.Range("A14:B21").Select()
